# Revert capacity chart to show kilowatts (kW) instead of watts on the
# y-axis: divide the raw capacity figures (currently expressed in Watts)
# by 1000, show one decimal place of precision in the cell number format,
# and relabel / reformat the chart's value axis accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Convert the underlying data from Watts to Kilowatts (divide by 1000)
#    Only the non-zero cells in columns C (Energy Storage), E (Solar) and
#    G (Wind) actually carried data.
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 5
$ws.Range("C25").Value = 7.6

$ws.Range("E10").Value = 7
$ws.Range("E12").Value = 15.04
$ws.Range("E13").Value = 111.11
$ws.Range("E14").Value = 3129.04
$ws.Range("E15").Value = 144.8
$ws.Range("E16").Value = 11
$ws.Range("E17").Value = 31.6
$ws.Range("E18").Value = 158.26
$ws.Range("E19").Value = 41.45
$ws.Range("E20").Value = 72.8
$ws.Range("E21").Value = 177.49
$ws.Range("E22").Value = 161.73
$ws.Range("E23").Value = 147.224
$ws.Range("E24").Value = 223.7
$ws.Range("E25").Value = 4271.835
$ws.Range("E26").Value = 1373.99

$ws.Range("G10").Value = 1.8
$ws.Range("G12").Value = 26

# ---------------------------------------------------------------------
# 2. The custom number format used by the data cells (numFmtId 164) goes
#    from "#,##0" to "#,##0.0" so the new fractional kilowatt values are
#    still displayed clearly. Re-apply the format string across the full
#    data range that used it (B2:G26); Excel will reuse a single style
#    for the whole range since they all share the same formatting.
# ---------------------------------------------------------------------
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# ---------------------------------------------------------------------
# 3. Update the chart's value (y) axis: title text and tick number format.
# ---------------------------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2)

$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"
